# Fix mistake in population sector where conveyors were initialized backwards.
# On the "Processing Initialization" sheet, the age-bucket counters in column B
# (1..20 per 20-row block) were counting up; they should count down (20..1)
# within each block, so reverse the two "master" mini age-tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Processing Initialization")

# Block 1: rows 2-20 (age counter 1..20 -> 20..1)
$ws.Range("B2").Value = 19
$ws.Range("B3").Formula = "=B2-1"
$ws.Range("B4:B20").Formula = "=B3-1"

# Block 2: rows 21-40 (age counter 1..20 -> 20..1)
$ws.Range("B21").Value = 20
$ws.Range("B22").Formula = "=B21-1"
$ws.Range("B23:B40").Formula = "=B22-1"

# Rows 41-100 simply mirror the B21/B22 counters (formulas unchanged,
# their cached values update automatically through recalculation).

# Make "Processing Initialization" the active sheet/tab with B41 selected.
$ws.Activate()
$ws.Range("B41").Select() | Out-Null
